$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3897.912874983024
$ws.Range("D2").Value = 267.6714347680194
